$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.859.19"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "1.764.27"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.96"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4463"
$ws.Range("E7").Value = "  -1.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3542"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07404"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.91"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.100"
$ws.Range("E11").Value = "  +2.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("E13").Value = "  +2.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.013"
$ws.Range("E14").Value = "  +2.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.229"
$ws.Range("E15").Value = "  +2.89%  "

$ws.Range("D16").Value = "1.761.76"
$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.90"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06422"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("E21").Value = "  +3.21%  "

$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("D23").Value = "27.882.55"

$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.109"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.35"
$ws.Range("E26").Value = "  -0.82%  "

$ws.Range("E27").Value = "  +1.93%  "

$ws.Range("D28").Value = "1.964.69"
$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.142"
$ws.Range("E29").Value = "  +5.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.88"
$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("E31").Value = "  +5.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09187"
$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.644"
$ws.Range("E33").Value = "  +5.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.689"
$ws.Range("E34").Value = "  +1.04%  "

$ws.Range("E35").Value = "  +2.40%  "

$ws.Range("E36").Value = "  +4.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02278"
$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2098"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6297"
$ws.Range("E39").Value = "  +1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.944"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.393"
$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.835"
$ws.Range("E43").Value = "  +2.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.24"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.742"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5840"
$ws.Range("E46").Value = "  +1.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.14"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.946"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("E50").Value = "  +2.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.64"
$ws.Range("E51").Value = "  +2.34%  "
